$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.647.80"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "1.534.27"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.89"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3933"
$ws.Range("E7").Value = "  +1.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3173"
$ws.Range("E8").Value = "  -2.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.40"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07174"
$ws.Range("E10").Value = "  -2.49%  "

$ws.Range("E11").Value = "  -6.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.0000"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.682"
$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("E14").Value = "  -3.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.589"
$ws.Range("E15").Value = "  -3.10%  "

$ws.Range("D16").Value = "1.540.91"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").Value = "  -2.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06587"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.63"
$ws.Range("E19").Value = "  -1.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.118"
$ws.Range("E21").Value = "  -4.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.44"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.69"
$ws.Range("E23").Value = "  -6.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.349"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("D25").Value = "21.655.93"
$ws.Range("E25").Value = "  -1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.350"
$ws.Range("E26").Value = "  -7.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.53"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("E28").Value = "  -2.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.844"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").Value = "1.723.91"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.03"
$ws.Range("E31").Value = "  -3.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.062"
$ws.Range("E32").Value = "  +4.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9392"
$ws.Range("E33").Value = "  -15.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08126"
$ws.Range("E34").Value = "  -1.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.514"
$ws.Range("E35").Value = "  -8.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.157"
$ws.Range("E36").Value = "  -1.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06017"
$ws.Range("E37").Value = "  -4.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02225"
$ws.Range("E38").Value = "  -3.31%  "

$ws.Range("E39").Value = "  -14.49%  "

$ws.Range("E40").Value = "  -4.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.177"
$ws.Range("E41").Value = "  -3.83%  "

$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9997"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5763"
$ws.Range("E44").Value = "  -3.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.10"
$ws.Range("E45").Value = "  -3.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.712"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5503"
$ws.Range("E47").Value = "  -4.41%  "

$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.878"
$ws.Range("E49").Value = "  -2.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.84"
$ws.Range("E50").Value = "  -2.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06691"
$ws.Range("E51").Value = "  -2.91%  "
